$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.760.52'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.51%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.216.25'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.47%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.17%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.18'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.05%  '

# Row 6
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.98%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.34'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.99%  '

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.18%  '

# Row 9
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.02%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.53%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0947'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.14%  '

# Row 12
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.47%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.94'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.02%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.551.74'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.32%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.23'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.70%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.834'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.80%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.207.45'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.62%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.667.76'

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.79%  '

# Row 20
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.43'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.34%  '

# Row 21
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.15'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.55%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.12'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +23.75%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.24'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.06%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -9.60%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.01%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.38'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.14%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.63'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.01%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.27'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.55%  '

# Row 29
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.62%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.38'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.03%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.39'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.37%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0797'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.09%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.50'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.03%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.97'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.49%  '

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.59%  '

# Row 36
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -11.53%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.27'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.47%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0302'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.24%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.32'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.55%  '

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.89%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.60'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.54%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '63.95'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.23%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.197'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.35%  '

# Row 44
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.07%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.76'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.52%  '

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.25%  '

# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.79%  '

# Row 48
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.11'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.33%  '

# Row 49
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.16'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.40%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.76%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.423.44'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.52%  '
